# Update national story with new data for 3 states (Ohio data fill-in
# across Admissions / Population sheets, and mark Ohio as "submitted"
# on the summary sheet).

$wb = $excel.ActiveWorkbook

# --- Admissions 2018: Ohio is row 36 ---
$ws = $wb.Worksheets.Item("Admissions 2018")
$ws.Range("C36").Value = 20697
$ws.Range("D36").Value = 6335
$ws.Range("E36").Value = 3234
$ws.Range("G36").Value = 3234
$ws.Range("H36").Value = 3101
$ws.Range("J36").Value = 3101

# --- Admissions 2019: Ohio is row 36 ---
$ws = $wb.Worksheets.Item("Admissions 2019")
$ws.Range("C36").Value = 20194
$ws.Range("D36").Value = 6357
$ws.Range("E36").Value = 3019
$ws.Range("G36").Value = 3019
$ws.Range("H36").Value = 3338
$ws.Range("J36").Value = 3338

# --- Admissions 2020: Ohio is row 36 ---
$ws = $wb.Worksheets.Item("Admissions 2020")
$ws.Range("C36").Value = 14022
$ws.Range("D36").Value = 4779
$ws.Range("E36").Value = 1931
$ws.Range("G36").Value = 1931
$ws.Range("H36").Value = 2848
$ws.Range("J36").Value = 2848
$ws.Range("K36").Value = 2020
$ws.Range("L36").Value = "CY"
$ws.Range("M36").Value = 12
$ws.Range("N36").Value = "Yes"

# --- Population 2018: Ohio is row 36 ---
$ws = $wb.Worksheets.Item("Population 2018")
$ws.Range("C36").Value = 48954
$ws.Range("D36").Value = 4342
$ws.Range("E36").Value = 2877
$ws.Range("G36").Value = 2877
$ws.Range("H36").Value = 1465
$ws.Range("J36").Value = 1465

# --- Population 2019: Ohio is row 36 ---
$ws = $wb.Worksheets.Item("Population 2019")
$ws.Range("C36").Value = 48697
$ws.Range("D36").Value = 4438
$ws.Range("E36").Value = 2870
$ws.Range("G36").Value = 2870
$ws.Range("H36").Value = 1568
$ws.Range("J36").Value = 1568

# --- Population 2020: Ohio is row 36 ---
$ws = $wb.Worksheets.Item("Population 2020")
$ws.Range("C36").Value = 43665
$ws.Range("D36").Value = 3662
$ws.Range("E36").Value = 2146
$ws.Range("G36").Value = 2146
$ws.Range("H36").Value = 1516
$ws.Range("J36").Value = 1516
$ws.Range("K36").Value = "Yes"

# --- States that Submitted: Ohio is row 36 ---
# Mark Ohio's admissions/population rows as submitted ("Yes") now that
# the data has been filled in, and clear the "needs data" highlight
# (orange fill) so the row matches the styling of the other completed
# rows (e.g. row 47).
$ws = $wb.Worksheets.Item("States that Submitted")
$ws.Range("B36").Value = "Yes"
$ws.Range("C36").Value = "Yes"

$ws.Range("A36").Interior.Color = $ws.Range("A47").Interior.Color
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Interior.Color = $ws.Range("B47").Interior.Color
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Interior.Color = $ws.Range("C47").Interior.Color
